# cps2017_dt_c.xlsx - "ajuste: corrigindo as categorias"
#
# Adds two new columns (S: "Idade ignorada", T: "Total") and two new rows
# ("Outros" and the grand-"Total" row) to the cid_grupos x age-group
# cross-tab on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new column headers -------------------------------------------------
$ws.Range("S1").Value = "Idade ignorada"
$ws.Range("T1").Value = "Total"

# ---- per-row totals for the existing disease categories (rows 2-6) ------
# Column S ("Idade ignorada") is left blank for these rows, matching the
# source data; column T gets the row total.
$ws.Range("T2").Value = 2199
$ws.Range("T3").Value = 277
$ws.Range("T4").Value = 1149
$ws.Range("T5").Value = 264
$ws.Range("T6").Value = 1438

# ---- row 7: "Outros" ------------------------------------------------------
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 141
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 70
$ws.Range("G7").Value = 66
$ws.Range("H7").Value = 76
$ws.Range("I7").Value = 95
$ws.Range("J7").Value = 84
$ws.Range("K7").Value = 103
$ws.Range("L7").Value = 112
$ws.Range("M7").Value = 96
$ws.Range("N7").Value = 133
$ws.Range("O7").Value = 135
$ws.Range("P7").Value = 116
$ws.Range("Q7").Value = 163
$ws.Range("R7").Value = 540
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 2000

# ---- row 8: grand "Total" row ---------------------------------------------
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 156
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 14
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 86
$ws.Range("G8").Value = 88
$ws.Range("H8").Value = 106
$ws.Range("I8").Value = 162
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 278
$ws.Range("L8").Value = 354
$ws.Range("M8").Value = 457
$ws.Range("N8").Value = 581
$ws.Range("O8").Value = 690
$ws.Range("P8").Value = 748
$ws.Range("Q8").Value = 813
$ws.Range("R8").Value = 2546
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 7327
